# Apply change: "From 1.4 to 1.5 version" - rotate the content blocks of
# TC1, TC2 and TC4 so that:
#   TC1 now describes "Competencias (portfolio)" (previously TC2's content)
#   TC2 now describes "Niveis das Competencias" (previously TC4's content)
#   TC4 now describes "Periodos Avaliativos" (previously TC1's content)
# TC3 (Avaliacoes) and TC5 (Perfis de Competencias) remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- TC1 block (rows 10-11) : Periodos Avaliativos -> Competencias (portfolio)
$ws.Range("B10").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Competencias (portfolio) a partir do menu inicial"
$ws.Range("D10").Value = "SYSTEM exibe a listagem das Competencias (portfolio) cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B11").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Competencias (portfolio)"

# --- TC2 block (rows 19-20) : Competencias (portfolio) -> Niveis das Competencias
$ws.Range("B19").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Niveis das Competencias a partir do menu inicial"
$ws.Range("D19").Value = "SYSTEM exibe a listagem dos Niveis das Competencias cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B20").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Niveis das Competencias"

# --- TC4 block (rows 37-38) : Niveis das Competencias -> Periodos Avaliativos
$ws.Range("B37").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Periodos Avaliativos a partir do menu inicial"
$ws.Range("D37").Value = "SYSTEM exibe a listagem dos Periodos Avaliativos cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B38").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Periodos Avaliativos"
